# Add team record (Wins / Losses / Ties) columns to the right of the
# existing data, mirroring the website_scraper shape but keeping W/L/T
# on this same sheet (AD:AF) instead of a separate sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/centered/bordered
# style used by the rest of row 1, then set the header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the same team record: 78 wins, 84 losses, 0 ties.
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 78   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 84   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
